$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark 3 corrected answers in red font color
$ws.Cells.Item(24, 9).Font.Color = 255
$ws.Cells.Item(25, 9).Font.Color = 255
$ws.Cells.Item(28, 8).Value = "+"
$ws.Cells.Item(28, 8).Font.Color = 255

# Legend: red-filled marker cell (H21) + "исправленно" label (I21)
$ws.Cells.Item(21, 8).Interior.Color = 255
$ws.Cells.Item(21, 9).Value = "исправленно"

# Reasons for the corrections in column K
$ws.Cells.Item(24, 11).Value = "Необходимо в случае использования Collections.sort"
$ws.Cells.Item(24, 11).VerticalAlignment = -4160

$ws.Cells.Item(25, 11).Value = "Необходимо в случае использования Collections.sort"
$ws.Cells.Item(25, 11).VerticalAlignment = -4160

$ws.Cells.Item(28, 11).Value = "Необходимо для определения уникальности элемента"
$ws.Cells.Item(28, 11).VerticalAlignment = -4160

$ws.Range("D43").Select() | Out-Null
